$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the numeric password value with the string "admin123"
$ws.Range("B2").Value = "admin123"

# Move the active selection to B3, matching the post-edit saved state
$ws.Range("B3").Select()
